# fix(publipostage): Correct status name
#
# Renames a handful of "statut" values in the shared-string table that are
# used throughout the worksheet:
#   - "bleu"                                                  -> "noir"
#   - "résultat et / ou publication posté"                    -> "résultat postés ou publiés"
#   - "pas de résultat ni de publication"                     -> "pas de résultat postés ni publiés"
#   - "résultat et / ou publication posté dans les 36 mois"   -> "résultat postés ou publiés dans les 36 mois"
#   - "résultat et / ou publication posté dans les 12 mois"   -> "résultat postés ou publiés dans les 12 mois"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of old text -> new text (exact, full-cell matches only)
$replacements = @{
    "bleu" = "noir";
    "résultat et / ou publication posté" = "résultat postés ou publiés";
    "pas de résultat ni de publication" = "pas de résultat postés ni publiés";
    "résultat et / ou publication posté dans les 36 mois" = "résultat postés ou publiés dans les 36 mois";
    "résultat et / ou publication posté dans les 12 mois" = "résultat postés ou publiés dans les 12 mois";
}

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $used.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -ne $null -and $replacements.ContainsKey($val)) {
            $cell.Value = $replacements[$val]
        }
    }
}
